$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -18.07537991065085
$ws.Cells.Item(2, 3).Value = -18.07537991065085
$ws.Cells.Item(2, 4).Value = -18.07537991065085
$ws.Cells.Item(2, 5).Value = -18.07537991065085
$ws.Cells.Item(2, 6).Value = -18.07537991065085
$ws.Cells.Item(2, 7).Value = -18.07537991065085
$ws.Cells.Item(2, 8).Value = -18.07537991065085
$ws.Cells.Item(2, 9).Value = -18.07537991065085
$ws.Cells.Item(2, 10).Value = -18.07537991065085
$ws.Cells.Item(2, 11).Value = -18.07537991065085
$ws.Cells.Item(3, 2).Value = -18.07537991065085
$ws.Cells.Item(3, 3).Value = -18.07537991065085
$ws.Cells.Item(3, 4).Value = -18.07537991065085
$ws.Cells.Item(3, 5).Value = -18.07537991065085
$ws.Cells.Item(3, 6).Value = -18.07537991065085
$ws.Cells.Item(3, 7).Value = -18.07537991065085
$ws.Cells.Item(3, 8).Value = -18.07537991065085
$ws.Cells.Item(3, 9).Value = 2.411986166535691
$ws.Cells.Item(3, 10).Value = -18.07537991065085
$ws.Cells.Item(3, 11).Value = -18.07537991065085
$ws.Cells.Item(4, 2).Value = -18.07537991065085
$ws.Cells.Item(4, 3).Value = -18.07537991065085
$ws.Cells.Item(4, 4).Value = 2.685357089462918
$ws.Cells.Item(4, 5).Value = -18.07537991065085
$ws.Cells.Item(4, 6).Value = 2.450011856680621
$ws.Cells.Item(4, 7).Value = -18.07537991065085
$ws.Cells.Item(4, 8).Value = 1.793896127331094
$ws.Cells.Item(4, 9).Value = -18.07537991065085
$ws.Cells.Item(4, 10).Value = 2.397506235378088
$ws.Cells.Item(4, 11).Value = -18.07537991065085
$ws.Cells.Item(5, 2).Value = -18.07537991065085
$ws.Cells.Item(5, 3).Value = 4.321923132766308
$ws.Cells.Item(5, 4).Value = -18.07537991065085
$ws.Cells.Item(5, 5).Value = -18.07537991065085
$ws.Cells.Item(5, 6).Value = -18.07537991065085
$ws.Cells.Item(5, 7).Value = 2.125989043006786
$ws.Cells.Item(5, 8).Value = -18.07537991065085
$ws.Cells.Item(5, 9).Value = -18.07537991065085
$ws.Cells.Item(5, 10).Value = -18.07537991065085
$ws.Cells.Item(5, 11).Value = -18.07537991065085
$ws.Cells.Item(6, 2).Value = -18.07537991065085
$ws.Cells.Item(6, 3).Value = -18.07537991065085
$ws.Cells.Item(6, 4).Value = -18.07537991065085
$ws.Cells.Item(6, 5).Value = -18.07537991065085
$ws.Cells.Item(6, 6).Value = -18.07537991065085
$ws.Cells.Item(6, 7).Value = -18.07537991065085
$ws.Cells.Item(6, 8).Value = -18.07537991065085
$ws.Cells.Item(6, 9).Value = -18.07537991065085
$ws.Cells.Item(6, 10).Value = -18.07537991065085
$ws.Cells.Item(6, 11).Value = -18.07537991065085
$ws.Cells.Item(7, 2).Value = 2.838246561821071
$ws.Cells.Item(7, 3).Value = -18.07537991065085
$ws.Cells.Item(7, 4).Value = -18.07537991065085
$ws.Cells.Item(7, 5).Value = -18.07537991065085
$ws.Cells.Item(7, 6).Value = -18.07537991065085
$ws.Cells.Item(7, 7).Value = -18.07537991065085
$ws.Cells.Item(7, 8).Value = -18.07537991065085
$ws.Cells.Item(7, 9).Value = -18.07537991065085
$ws.Cells.Item(7, 10).Value = -18.07537991065085
$ws.Cells.Item(7, 11).Value = -18.07537991065085
$ws.Cells.Item(8, 2).Value = -18.07537991065085
$ws.Cells.Item(8, 3).Value = -18.07537991065085
$ws.Cells.Item(8, 4).Value = -18.07537991065085
$ws.Cells.Item(8, 5).Value = 2.378607045338765
$ws.Cells.Item(8, 6).Value = -18.07537991065085
$ws.Cells.Item(8, 7).Value = -18.07537991065085
$ws.Cells.Item(8, 8).Value = -18.07537991065085
$ws.Cells.Item(8, 9).Value = -18.07537991065085
$ws.Cells.Item(8, 10).Value = -18.07537991065085
$ws.Cells.Item(8, 11).Value = -18.07537991065085
$ws.Cells.Item(9, 2).Value = 3.683520335815912
$ws.Cells.Item(9, 3).Value = -18.07537991065085
$ws.Cells.Item(9, 4).Value = -18.07537991065085
$ws.Cells.Item(9, 5).Value = -18.07537991065085
$ws.Cells.Item(9, 6).Value = -18.07537991065085
$ws.Cells.Item(9, 7).Value = -18.07537991065085
$ws.Cells.Item(9, 8).Value = -18.07537991065085
$ws.Cells.Item(9, 9).Value = -18.07537991065085
$ws.Cells.Item(9, 10).Value = -18.07537991065085
$ws.Cells.Item(9, 11).Value = -18.07537991065085
$ws.Cells.Item(10, 2).Value = -18.07537991065085
$ws.Cells.Item(10, 3).Value = -18.07537991065085
$ws.Cells.Item(10, 4).Value = -18.07537991065085
$ws.Cells.Item(10, 5).Value = -18.07537991065085
$ws.Cells.Item(10, 6).Value = -18.07537991065085
$ws.Cells.Item(10, 7).Value = -18.07537991065085
$ws.Cells.Item(10, 8).Value = -18.07537991065085
$ws.Cells.Item(10, 9).Value = 1.625420449717309
$ws.Cells.Item(10, 10).Value = -18.07537991065085
$ws.Cells.Item(10, 11).Value = 2.197788497394426
$ws.Cells.Item(11, 2).Value = -18.07537991065085
$ws.Cells.Item(11, 3).Value = -18.07537991065085
$ws.Cells.Item(11, 4).Value = -18.07537991065085
$ws.Cells.Item(11, 5).Value = 2.208247696839257
$ws.Cells.Item(11, 6).Value = -18.07537991065085
$ws.Cells.Item(11, 7).Value = 2.575364211962484
$ws.Cells.Item(11, 8).Value = -18.07537991065085
$ws.Cells.Item(11, 9).Value = -18.07537991065085
$ws.Cells.Item(11, 10).Value = -18.07537991065085
$ws.Cells.Item(11, 11).Value = 1.344170913516387
$ws.Cells.Item(12, 2).Value = -18.07537991065085
$ws.Cells.Item(12, 3).Value = -18.07537991065085
$ws.Cells.Item(12, 4).Value = -18.07537991065085
$ws.Cells.Item(12, 5).Value = -18.07537991065085
$ws.Cells.Item(12, 6).Value = -18.07537991065085
$ws.Cells.Item(12, 7).Value = -18.07537991065085
$ws.Cells.Item(12, 8).Value = -18.07537991065085
$ws.Cells.Item(12, 9).Value = -18.07537991065085
$ws.Cells.Item(12, 10).Value = -18.07537991065085
$ws.Cells.Item(12, 11).Value = -18.07537991065085
$ws.Cells.Item(13, 2).Value = -18.07537991065085
$ws.Cells.Item(13, 3).Value = -18.07537991065085
$ws.Cells.Item(13, 4).Value = -18.07537991065085
$ws.Cells.Item(13, 5).Value = 1.962216641535255
$ws.Cells.Item(13, 6).Value = -18.07537991065085
$ws.Cells.Item(13, 7).Value = -18.07537991065085
$ws.Cells.Item(13, 8).Value = -18.07537991065085
$ws.Cells.Item(13, 9).Value = -18.07537991065085
$ws.Cells.Item(13, 10).Value = 2.307561738968592
$ws.Cells.Item(13, 11).Value = 1.606608222708972
$ws.Cells.Item(14, 2).Value = -18.07537991065085
$ws.Cells.Item(14, 3).Value = -18.07537991065085
$ws.Cells.Item(14, 4).Value = 1.941315498878935
$ws.Cells.Item(14, 5).Value = -18.07537991065085
$ws.Cells.Item(14, 6).Value = -18.07537991065085
$ws.Cells.Item(14, 7).Value = -18.07537991065085
$ws.Cells.Item(14, 8).Value = -18.07537991065085
$ws.Cells.Item(14, 9).Value = -18.07537991065085
$ws.Cells.Item(14, 10).Value = -18.07537991065085
$ws.Cells.Item(14, 11).Value = 2.070090824078662
$ws.Cells.Item(15, 2).Value = -18.07537991065085
$ws.Cells.Item(15, 3).Value = -18.07537991065085
$ws.Cells.Item(15, 4).Value = 0.7858918045702903
$ws.Cells.Item(15, 5).Value = -18.07537991065085
$ws.Cells.Item(15, 6).Value = -18.07537991065085
$ws.Cells.Item(15, 7).Value = -18.07537991065085
$ws.Cells.Item(15, 8).Value = -18.07537991065085
$ws.Cells.Item(15, 9).Value = -18.07537991065085
$ws.Cells.Item(15, 10).Value = -18.07537991065085
$ws.Cells.Item(15, 11).Value = -18.07537991065085
$ws.Cells.Item(16, 2).Value = -18.07537991065085
$ws.Cells.Item(16, 3).Value = -18.07537991065085
$ws.Cells.Item(16, 4).Value = -18.07537991065085
$ws.Cells.Item(16, 5).Value = -18.07537991065085
$ws.Cells.Item(16, 6).Value = -18.07537991065085
$ws.Cells.Item(16, 7).Value = -18.07537991065085
$ws.Cells.Item(16, 8).Value = -18.07537991065085
$ws.Cells.Item(16, 9).Value = -18.07537991065085
$ws.Cells.Item(16, 10).Value = 2.33318306985037
$ws.Cells.Item(16, 11).Value = -18.07537991065085
$ws.Cells.Item(17, 2).Value = -18.07537991065085
$ws.Cells.Item(17, 3).Value = -18.07537991065085
$ws.Cells.Item(17, 4).Value = -0.1056505201517785
$ws.Cells.Item(17, 5).Value = -18.07537991065085
$ws.Cells.Item(17, 6).Value = -18.07537991065085
$ws.Cells.Item(17, 7).Value = -18.07537991065085
$ws.Cells.Item(17, 8).Value = 0.4143622328045399
$ws.Cells.Item(17, 9).Value = 0.9743649626346429
$ws.Cells.Item(17, 10).Value = 1.214947236911321
$ws.Cells.Item(17, 11).Value = -18.07537991065085
$ws.Cells.Item(18, 2).Value = -18.07537991065085
$ws.Cells.Item(18, 3).Value = -18.07537991065085
$ws.Cells.Item(18, 4).Value = -18.07537991065085
$ws.Cells.Item(18, 5).Value = -18.07537991065085
$ws.Cells.Item(18, 6).Value = -18.07537991065085
$ws.Cells.Item(18, 7).Value = -18.07537991065085
$ws.Cells.Item(18, 8).Value = 0.3117160503306202
$ws.Cells.Item(18, 9).Value = 0.8127836184849212
$ws.Cells.Item(18, 10).Value = 1.275060790127294
$ws.Cells.Item(18, 11).Value = -18.07537991065085
$ws.Cells.Item(19, 2).Value = -18.07537991065085
$ws.Cells.Item(19, 3).Value = -18.07537991065085
$ws.Cells.Item(19, 4).Value = 1.189385665363067
$ws.Cells.Item(19, 5).Value = -18.07537991065085
$ws.Cells.Item(19, 6).Value = -18.07537991065085
$ws.Cells.Item(19, 7).Value = -18.07537991065085
$ws.Cells.Item(19, 8).Value = 1.886568529325788
$ws.Cells.Item(19, 9).Value = 2.082427120287689
$ws.Cells.Item(19, 10).Value = -18.07537991065085
$ws.Cells.Item(19, 11).Value = -18.07537991065085
$ws.Cells.Item(20, 2).Value = -18.07537991065085
$ws.Cells.Item(20, 3).Value = -18.07537991065085
$ws.Cells.Item(20, 4).Value = 2.260896720493216
$ws.Cells.Item(20, 5).Value = -18.07537991065085
$ws.Cells.Item(20, 6).Value = 3.86153151240503
$ws.Cells.Item(20, 7).Value = -18.07537991065085
$ws.Cells.Item(20, 8).Value = 2.232546840415349
$ws.Cells.Item(20, 9).Value = 1.862351653644951
$ws.Cells.Item(20, 10).Value = -18.07537991065085
$ws.Cells.Item(20, 11).Value = 2.492844365697804
$ws.Cells.Item(21, 2).Value = -18.07537991065085
$ws.Cells.Item(21, 3).Value = -18.07537991065085
$ws.Cells.Item(21, 4).Value = -18.07537991065085
$ws.Cells.Item(21, 5).Value = 2.651191673390736
$ws.Cells.Item(21, 6).Value = -18.07537991065085
$ws.Cells.Item(21, 7).Value = 3.274217188490111
$ws.Cells.Item(21, 8).Value = 2.475442952347751
$ws.Cells.Item(21, 9).Value = -18.07537991065085
$ws.Cells.Item(21, 10).Value = -18.07537991065085
$ws.Cells.Item(21, 11).Value = -18.07537991065085
